# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Pais")

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 22:05"

# --- Updated totals for existing countries (no re-ranking) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1422217
$ws.Range("C4").Value = 13581
$ws.Range("D4").Value = 306298
$ws.Range("E4").Value = 1031365
$ws.Range("G4").Value = 1129
$ws.Range("H4").Value = 84554

# Row 11: Alemania
$ws.Range("B11").Value = 173919
$ws.Range("C11").Value = 748
$ws.Range("E11").Value = 17400
$ws.Range("G11").Value = 81
$ws.Range("H11").Value = 7819

# Row 24: Ecuador
$ws.Range("B24").Value = 30486
$ws.Range("C24").Value = 67
$ws.Range("E24").Value = 24719
$ws.Range("F24").Value = 190
$ws.Range("G24").Value = 7
$ws.Range("H24").Value = 2334

# Row 25: Suiza
$ws.Range("D25").Value = 27100
$ws.Range("E25").Value = 1443

# --- Togo moves above Liberia/Madagascar in the ranking (rows 145-147) ---
# Row 145 becomes Togo with its new figures
$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 219
$ws.Range("C145").Value = 20
$ws.Range("D145").Value = 96
$ws.Range("E145").Value = 112
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 11

# Row 146 becomes Liberia (previous row-145 figures)
$ws.Range("A146").Value = "Liberia"
$ws.Range("B146").Value = 213
$ws.Range("C146").Value = 2
$ws.Range("D146").Value = 101
$ws.Range("E146").Value = 92
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 20

# Row 147 becomes Madagascar (previous row-146 figures)
$ws.Range("A147").Value = "Madagascar"
$ws.Range("B147").Value = 212
$ws.Range("C147").Value = 26
$ws.Range("D147").Value = 107
$ws.Range("E147").Value = 105
$ws.Range("F147").Value = 1
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0

# Row 148 (Sudan del Sur) is unaffected by the reordering

# --- Belice moves above Nueva Caledonia in the ranking (rows 193-194) ---
# Row 193 becomes Belice
$ws.Range("A193").Value = "Belice"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 16
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 2

# Row 194 becomes Nueva Caledonia
$ws.Range("A194").Value = "Nueva Caledonia"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 18
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0
